# Add eight more rows (5-12) of Amount_Spend / Business_Received / Deals_Closed
# data to Sheet1, give them a left/right "column separator" border + centered
# alignment, and move the active selection to D12 (one cell below/right of the
# newly added table), matching the author's "Add files via upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- new data rows --------------------------------------------------------
$data = @(
    @(300,  50,  42),
    @(100,  12,  10),
    @(1600, 120, 100),
    @(3100, 300, 200),
    @(4600, 600, 200),
    @(6100, 100, 50),
    @(7600, 178, 78),
    @(9100, 246, 46)
)

$startRow = 5
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 1; $c -le 3; $c++) {
        $ws.Cells.Item($row, $c).Value = $vals[$c - 1]
    }
}

# ---- formatting -------------------------------------------------------
# Build the new "thin left/right border + centered" look on the first new
# cell, then fan it out to the rest of the block with a format-only paste so
# the whole A5:C12 block shares a single new style (rather than minting a
# fresh style per cell).
$seed = $ws.Cells.Item($startRow, 1)
$seed.HorizontalAlignment = -4108   # xlCenter
$seed.VerticalAlignment = -4108     # xlCenter
$seed.Borders.Item(7).LineStyle = 1   # xlEdgeLeft  -> xlContinuous
$seed.Borders.Item(10).LineStyle = 1  # xlEdgeRight -> xlContinuous

$seed.Copy() | Out-Null
$ws.Range("A5:C12").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ---- selection / cursor -----------------------------------------------
$ws.Range("D12").Select() | Out-Null
